$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Value = 39
$ws.Range("B40").Value = 1
$ws.Range("C40").Value = "2024-06-15 17:11:47"
$ws.Range("D40").Value = 200
$ws.Range("E40").Value = 12

$ws.Range("A41").Value = 40
$ws.Range("B41").Value = 2
$ws.Range("C41").Value = "2024-06-15 17:11:48"
$ws.Range("D41").Value = 200
$ws.Range("E41").Value = 0
